# Auto commit at 2026-02-09 10:20:58.69
# Updates daydata.xlsx: corrects the 2024-01-01 (46054) "四方坪站" row,
# tags J10 with a new number-format style, and appends four new daily
# rows (46060 / 46061) for both stations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct row 2 (46054, 四方坪站) figures ---------------------------
$ws.Range("C2").Value = 11580.25
$ws.Range("D2").Value = 10486.99
$ws.Range("E2").Value = 4302.8599999999997
$ws.Range("F2").Value = 469

# --- Tag J10 with a new "#,##0.00" (built-in numFmtId 4) style -------
$ws.Range("J10").NumberFormat = "#,##0.00"

# Give column J a sensible width like the rest of the data columns.
# (ColumnWidth is expressed in "characters"; the engine stores the
# OOXML <col width> as characters + 5/7, quantized to 1/7ths, so this
# is the closest achievable value to the target 12.125.)
$ws.Columns.Item(10).ColumnWidth = 11.428571428571429

# --- New row 14: 2024-xx-xx (46060) 四方坪站 ---------------------------
$ws.Range("A14").Value = 46060
$ws.Range("B14").Value = "四方坪站"
$ws.Range("C14").Value = 11441.41
$ws.Range("D14").Value = 10386.73
$ws.Range("E14").Value = 4227.2299999999996
$ws.Range("F14").Value = 471

# --- New row 15: 46060 高岭站 -------------------------------------------
$ws.Range("A15").Value = 46060
$ws.Range("B15").Value = "高岭站"
$ws.Range("C15").Value = 4838.78
$ws.Range("D15").Value = 4421.3
$ws.Range("E15").Value = 1317.65
$ws.Range("F15").Value = 180

# --- New row 16: 46061 四方坪站 -----------------------------------------
$ws.Range("A16").Value = 46061
$ws.Range("B16").Value = "四方坪站"
$ws.Range("C16").Value = 11078.34
$ws.Range("D16").Value = 9884.91
$ws.Range("E16").Value = 4197.91
$ws.Range("F16").Value = 452

# --- New row 17: 46061 高岭站 -------------------------------------------
$ws.Range("A17").Value = 46061
$ws.Range("B17").Value = "高岭站"
$ws.Range("C17").Value = 4649.2
$ws.Range("D17").Value = 4059.99
$ws.Range("E17").Value = 1246.24
$ws.Range("F17").Value = 161

# --- Scroll the view down so row 4 is the first visible row ----------
# (best-effort; mirrors the author's <sheetView topLeftCell="A4">)
$ws.Range("H10").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
